# Apply the trainer-availability restructuring edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Dency Patel" row (row 2); this shifts "Pratyush Singh" up to row 2.
$ws.Rows.Item(2).Delete()

# Split the single "AVAILABILITY LIST" column into DATE / FROM TIME / TO TIME columns,
# keeping only the final availability entry (2023-04-29).
$ws.Range("E1").Value = "DATE"
$ws.Range("F1").Value = "FROM TIME"
$ws.Range("G1").Value = "TO TIME"

# Force E2 to stay text (otherwise Excel auto-converts the date-looking string to a date serial).
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2023-04-29"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "09:00:00"
$ws.Range("G2").Value = "17:00:00"
